$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (BonusPower odds buffed from 1/50 to 1/45,
# and monsters buffed +100 starting power).
$rows = @(
    @("BonusPower", 2,    760, 96,  100, 200, "win"),
    @("BonusPower", 2,    770, 97,  10,  20,  "win"),
    @("BonusPower", 2,    570, 98,  100, 200, "win"),
    @("SkipBoss",   2,    759, 104, 20,  40,  "win"),
    @("BonusPower", 0.06, 0,   14,  100, 0,   "lose"),
    @("BonusPower", 2,    780, 98,  100, 200, "win")
)

$startRow = 17
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}
